# Refresh the "cryptos" price/volume snapshot (GitHub Actions data pull).
# For numeric-looking "Price" values, force the cell to Text format first so
# Excel's type inference doesn't silently turn the literal string into a
# number (the sheet stores these as plain text, e.g. "490.97", not 490.97).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.043.62'
$ws.Range('E2').Value = '  +4.17%  '
$ws.Range('D3').Value = '2.516.76'
$ws.Range('E3').Value = '  +5.12%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '490.97'
$ws.Range('E5').Value = '  +6.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.54'
$ws.Range('E6').Value = '  +12.05%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.515'
$ws.Range('E8').Value = '  +6.71%  '
$ws.Range('D9').Value = '2.531.40'
$ws.Range('E9').Value = '  +5.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0987'
$ws.Range('E10').Value = '  +5.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.62'
$ws.Range('E11').Value = '  +5.08%  '
$ws.Range('E12').Value = '  +5.74%  '
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '2.945.36'
$ws.Range('E14').Value = '  +5.31%  '
$ws.Range('D15').Value = '56.025.90'
$ws.Range('E15').Value = '  +4.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.07'
$ws.Range('E16').Value = '  +8.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000136'
$ws.Range('E17').Value = '  +8.13%  '
$ws.Range('D18').Value = '2.519.44'
$ws.Range('E18').Value = '  +5.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.45'
$ws.Range('E19').Value = '  +6.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.25'
$ws.Range('E20').Value = '  +10.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.24'
$ws.Range('E21').Value = '  +4.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -1.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.84'
$ws.Range('E23').Value = '  +9.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.58'
$ws.Range('E24').Value = '  +5.13%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.412'
$ws.Range('E25').Value = '  +8.74%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  +8.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = '2.613.80'
$ws.Range('E28').Value = '  +6.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.52'
$ws.Range('E29').Value = '  +6.00%  '
$ws.Range('D30').Value = '0.0₃0792'
$ws.Range('E30').Value = '  +12.73%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.28'
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.48'
$ws.Range('E33').Value = '  +4.91%  '
$ws.Range('E34').Value = '  +8.75%  '
$ws.Range('E35').Value = '  +5.62%  '
$ws.Range('E36').Value = '  +12.27%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.875'
$ws.Range('E37').Value = '  +11.15%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.71'
$ws.Range('E38').Value = '  +7.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.30'
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('E40').Value = '  +8.15%  '
$ws.Range('E41').Value = '  +3.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0559'
$ws.Range('E42').Value = '  +6.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.992'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.33'
$ws.Range('E44').Value = '  +9.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '267.78'
$ws.Range('E45').Value = '  +30.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.87'
$ws.Range('E46').Value = '  +13.26%  '
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0227'
$ws.Range('E48').Value = '  +6.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0909'
$ws.Range('E49').Value = '  +6.51%  '
$ws.Range('D50').Value = '1.958.38'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.89'
$ws.Range('E51').Value = '  +8.74%  '
